$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $value)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.ClearFormats()
}

Set-TextCell $ws 'D2' '26.077.87'
Set-TextCell $ws 'E2' '  -0.32%  '
Set-TextCell $ws 'D3' '1.645.89'
Set-TextCell $ws 'E3' '  -1.43%  '
Set-TextCell $ws 'E4' '  -0.18%  '
Set-TextCell $ws 'D5' '215.95'
Set-TextCell $ws 'E5' '  +2.51%  '
Set-TextCell $ws 'D6' '0.5218'
Set-TextCell $ws 'E6' '  +0.35%  '
Set-TextCell $ws 'E7' '  -0.14%  '
Set-TextCell $ws 'D8' '0.2613'
Set-TextCell $ws 'E8' '  -0.45%  '
Set-TextCell $ws 'D9' '0.06361'
Set-TextCell $ws 'E9' '  +0.62%  '
Set-TextCell $ws 'D10' '20.86'
Set-TextCell $ws 'E10' '  -1.66%  '
Set-TextCell $ws 'D11' '0.07673'
Set-TextCell $ws 'E11' '  +1.61%  '
Set-TextCell $ws 'D12' '1.650.05'
Set-TextCell $ws 'E12' '  -1.24%  '
Set-TextCell $ws 'D13' '4.420'
Set-TextCell $ws 'E13' '  -0.46%  '
Set-TextCell $ws 'D14' '1.868.73'
Set-TextCell $ws 'D15' '0.5548'
Set-TextCell $ws 'E15' '  +1.00%  '
Set-TextCell $ws 'D16' '0.0₅8264'
Set-TextCell $ws 'E16' '  +3.09%  '
Set-TextCell $ws 'D17' '64.98'
Set-TextCell $ws 'E17' '  -2.19%  '
Set-TextCell $ws 'D18' '26.102.31'
Set-TextCell $ws 'E18' '  -0.23%  '
Set-TextCell $ws 'D20' '4.729'
Set-TextCell $ws 'E20' '  -0.65%  '
Set-TextCell $ws 'D21' '188.26'
Set-TextCell $ws 'E21' '  +0.62%  '
Set-TextCell $ws 'D22' '10.22'
Set-TextCell $ws 'E22' '  -1.16%  '
Set-TextCell $ws 'E23' '  +0.29%  '
Set-TextCell $ws 'E24' '  -0.13%  '
Set-TextCell $ws 'D25' '146.14'
Set-TextCell $ws 'E25' '  -2.62%  '
Set-TextCell $ws 'D26' '0.1219'
Set-TextCell $ws 'E26' '  -1.68%  '
Set-TextCell $ws 'D27' '7.426'
Set-TextCell $ws 'E27' '  -0.98%  '
Set-TextCell $ws 'D28' '15.83'
Set-TextCell $ws 'E28' '  +0.02%  '
Set-TextCell $ws 'D29' '1.380'
Set-TextCell $ws 'E29' '  +2.19%  '
Set-TextCell $ws 'D30' '0.05979'
Set-TextCell $ws 'E30' '  -5.29%  '
Set-TextCell $ws 'D31' '1.269'
Set-TextCell $ws 'E31' '  -1.09%  '
Set-TextCell $ws 'B32' 'Filecoin'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D32' '3.398'
Set-TextCell $ws 'E32' '  -0.48%  '
Set-TextCell $ws 'B33' 'InternetComputer(DFINITY)'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D33' '3.406'
Set-TextCell $ws 'E33' '  -3.34%  '
Set-TextCell $ws 'D34' '1.663'
Set-TextCell $ws 'E34' '  +1.12%  '
Set-TextCell $ws 'D35' '0.9940'
Set-TextCell $ws 'E35' '  -1.11%  '
Set-TextCell $ws 'E36' '  -0.49%  '
Set-TextCell $ws 'D37' '2.749'
Set-TextCell $ws 'E37' '  -0.19%  '
Set-TextCell $ws 'D38' '0.5651'
Set-TextCell $ws 'E38' '  -6.79%  '
Set-TextCell $ws 'E39' '  +0.26%  '
Set-TextCell $ws 'D40' '0.8611'
Set-TextCell $ws 'E40' '  -0.29%  '
Set-TextCell $ws 'D41' '5.840'
Set-TextCell $ws 'D42' '1.002'
Set-TextCell $ws 'E42' '  -0.27%  '
Set-TextCell $ws 'D43' '1.029.91'
Set-TextCell $ws 'E43' '  -7.46%  '
Set-TextCell $ws 'D44' '100.25'
Set-TextCell $ws 'E44' '  -0.25%  '
Set-TextCell $ws 'D45' '1.795.89'
Set-TextCell $ws 'E45' '  -1.45%  '
Set-TextCell $ws 'E46' '  -1.22%  '
Set-TextCell $ws 'D47' '55.85'
Set-TextCell $ws 'E47' '  +0.53%  '
Set-TextCell $ws 'D48' '1.002'
Set-TextCell $ws 'E48' '  +0.25%  '
Set-TextCell $ws 'D49' '8.101'
Set-TextCell $ws 'E49' '  +0.65%  '
Set-TextCell $ws 'D50' '0.05167'
Set-TextCell $ws 'E50' '  -1.34%  '
Set-TextCell $ws 'D51' '0.4222'
Set-TextCell $ws 'E51' '  -0.46%  '
